# Cost Estimation workbook update:
#  - Copy "Kostenabschätzung Serie Detail" to a new sheet "Kostenabschätzung Serie Det (2)"
#    at the end of the workbook, and update a handful of its input figures
#    (new purchase-price numbers for the series run).
#  - Rename "Kostenübersicht Prototypen" to "Prototyp 1".
#  - Refresh the selection/active-cell bookmarks left on a few sheets.

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate the "Kostenabschätzung Serie Detail" sheet -------------
$detailSheet = $wb.Worksheets.Item("Kostenabschätzung Serie Detail")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$detailSheet.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Kostenabschätzung Serie Det (2)"

# --- 2. Rename "Kostenübersicht Prototypen" -> "Prototyp 1" --------------
# (this also repoints the formula on "Kostenabschätzung"!B5 automatically)
$protoSheet = $wb.Worksheets.Item("Kostenübersicht Prototypen")
$protoSheet.Name = "Prototyp 1"

# --- 3. Update the figures on the new (2) sheet ---------------------------
# N3 compares against its own sheet, so the copied formula must be repointed
# from the source sheet name to this sheet's own (new) name. Using Replace
# (rather than re-typing the whole formula) keeps the cell's original
# (unformatted) number style intact; re-paste the original cell's format
# afterwards to be safe.
$n3 = $newSheet.Range("N3")
$n3.Replace("Kostenabschätzung Serie Detail", "Kostenabschätzung Serie Det (2)", 2) | Out-Null
$detailSheet.Range("N3").Copy() | Out-Null
$n3.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$newSheet.Range("C8").Value = 507.35
$newSheet.Range("D8").Value = 2290.94
$newSheet.Range("E8").Value = 2290.94

$newSheet.Range("C18").Value = 139.38
$newSheet.Range("D18").Value = 655.38
$newSheet.Range("E18").Value = 655.38

$newSheet.Range("C20").Value = 0
$newSheet.Range("C21").Value = 0
$newSheet.Range("C22").Value = 0

# --- 4. Refresh selections on the touched sheets --------------------------
$protoSheet.Activate()
$protoSheet.Range("E33").Select()

$kostenSheet = $wb.Worksheets.Item("Kostenabschätzung")
$kostenSheet.Activate()
$kostenSheet.Range("B5").Select()

$detailSheet.Activate()
$detailSheet.Range("D8").Select()

$newSheet.Activate()
$newSheet.Range("D5").Select()

Write-Output "done"
